$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.968.36"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "2.244.24"
$ws.Range("E3").Value = "  -1.23%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "306.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "96.22"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.522"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("E8").Value = "  +0.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.488"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "34.63"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("E11").Value = "  +2.35%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.28%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.75"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "2.590.36"
$ws.Range("E14").Value = "  -1.32%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "14.37"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "2.245.46"
$ws.Range("E16").Value = "  -0.65%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.780"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "41.805.86"
$ws.Range("E18").Value = "  -0.97%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.14"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -1.20%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.91"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "67.03"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "235.22"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.24%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.56"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -0.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("B27").Value = "InjectiveProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "37.65"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "23.15"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("E29").Value = "  +1.65%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.44"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "166.32"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.45%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.09%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.15"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.01%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.04"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.55%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "17.34"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.0718"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -0.09%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.102"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("E40").Value = "  -2.71%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "1.938.91"
$ws.Range("E42").Value = "  -2.89%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0280"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.75%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -10.81%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "18.39"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "9.60"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "53.76"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "2.462.87"
$ws.Range("E49").Value = "  -1.24%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "70.98"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "90.89"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "
